# Add 6 new "Computer Part" accessory rows (rows 60-65) to Sheet1,
# matching the appended shared-string entries and row heights from the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$typeVal = 'Computer Part'
$noneVal = 'none'
$brandVal = 'Lenovo'

# Row 60: 'ThinkPad Battery 77++'
$ws.Cells.Item(60, 1).Value = 58
$ws.Cells.Item(60, 2).Value = $typeVal
$ws.Cells.Item(60, 3).Value = 'ThinkPad Battery 77++'
$ws.Cells.Item(60, 4).Value = $noneVal
$ws.Cells.Item(60, 5).Value = $noneVal
$ws.Cells.Item(60, 6).Value = $noneVal
$ws.Cells.Item(60, 7).Value = 140
$ws.Cells.Item(60, 8).Value = '•  The 6-cell ThinkPad battery 77++ is a replacement/spare battery for 
     ThinkPad P52. 
•  This battery is powered by Lithium-Ion technology and has a safeguard 
     chip which authenticates with Lenovo systems.
•  The battery lasts on average 90 watt hours and utilizes over-discharge
     protection to ensure a secure flow of power.
•  6 cell battery
•  Over-discharge protection to prolong battery life
•  Battery energy (Watt-hours): 90Wh
•  Battery safeguard chip authenticates with Lenovo systems  
•  Battery Life:  90Wh
•  Hardware Requirements:  Nominal voltage 11.4V'
$ws.Cells.Item(60, 8).WrapText = $true
$ws.Cells.Item(60, 9).Value = $brandVal
$ws.Cells.Item(60, 10).Value = 0
$ws.Cells.Item(60, 11).Value = 0
$ws.Rows.Item(60).RowHeight = 290

# Row 61: 'ThinkPad Battery 61++'
$ws.Cells.Item(61, 1).Value = 59
$ws.Cells.Item(61, 2).Value = $typeVal
$ws.Cells.Item(61, 3).Value = 'ThinkPad Battery 61++'
$ws.Cells.Item(61, 4).Value = $noneVal
$ws.Cells.Item(61, 5).Value = $noneVal
$ws.Cells.Item(61, 6).Value = $noneVal
$ws.Cells.Item(61, 7).Value = 140
$ws.Cells.Item(61, 8).Value = '• The 6-cell ThinkPad battery 61++ is a replacement/spare battery 
    for ThinkPad new T-series system. 
•  The battery is rechargeable and can be used as replacement or 
    as a convenient spare for ThinkPad T470 & T570 systems. 
•  This battery is powered by Lithium-Ion technology and has a 
safeguard chip which authenticates with Lenovo systems.
• Lithium-ion technology
• 6 cell cylindrical flat battery
• Over-discharge protection to prolong the life of the battery
• Battery energy (Watt-hours): 72Wh
• Battery Safeguard chip authenticates with Lenovo systems"
•  Battery Life:  5.5Hours 	'
$ws.Cells.Item(61, 8).WrapText = $true
$ws.Cells.Item(61, 9).Value = $brandVal
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 0
$ws.Rows.Item(61).RowHeight = 275.5

# Row 62: 'ThinkPad 500G Hard Drive'
$ws.Cells.Item(62, 1).Value = 60
$ws.Cells.Item(62, 2).Value = $typeVal
$ws.Cells.Item(62, 3).Value = 'ThinkPad 500G Hard Drive'
$ws.Cells.Item(62, 4).Value = $noneVal
$ws.Cells.Item(62, 5).Value = $noneVal
$ws.Cells.Item(62, 6).Value = $noneVal
$ws.Cells.Item(62, 7).Value = 90
$ws.Cells.Item(62, 8).Value = '• The ThinkPad 500GB 7200rpm 7mm 2.5” hard drive with the high
   speed SATA 6Gbps interface, 7mm height, it is the best solution
   to upgrade your ThinkPad systems.
• Provides a robust thin design enabling a seamless upgrade from 
   the primary hard drive of supported ThinkPad laptops.
• Provides good storage size along with advanced format and SATA
   data transfer capabilities.
• Hard Drive Capacity:  500GB
• Hard Drive Size:  2.5in'
$ws.Cells.Item(62, 8).WrapText = $true
$ws.Cells.Item(62, 9).Value = $brandVal
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Rows.Item(62).RowHeight = 203

# Row 63: 'ThinkStation 1TB Hard Drive'
$ws.Cells.Item(63, 1).Value = 61
$ws.Cells.Item(63, 2).Value = $typeVal
$ws.Cells.Item(63, 3).Value = 'ThinkStation 1TB Hard Drive'
$ws.Cells.Item(63, 4).Value = $noneVal
$ws.Cells.Item(63, 5).Value = $noneVal
$ws.Cells.Item(63, 6).Value = $noneVal
$ws.Cells.Item(63, 7).Value = 270
$ws.Cells.Item(63, 8).Value = '• ThinkStation 1TB 7200rpm SATA 3.5" Enterprise Hard Drive have
   high reliability and large capacity to support 7x24 working environment,
   it is fully compatible with specified ThinkStation machines, it is the best 
   solution for you to upgrade your ThinkStation Storage.
• OS Requirements: All supported. 
• Hardware Requirements: SATA PN
• Connection Type: SATA
•  Weight: 716g/1.58lb'
$ws.Cells.Item(63, 8).WrapText = $true
$ws.Cells.Item(63, 9).Value = $brandVal
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Rows.Item(63).RowHeight = 174

# Row 64: 'Lenovo 6GB DDR4 RAM Memory'
$ws.Cells.Item(64, 1).Value = 62
$ws.Cells.Item(64, 2).Value = $typeVal
$ws.Cells.Item(64, 3).Value = 'Lenovo 6GB DDR4 RAM Memory'
$ws.Cells.Item(64, 4).Value = $noneVal
$ws.Cells.Item(64, 5).Value = $noneVal
$ws.Cells.Item(64, 6).Value = $noneVal
$ws.Cells.Item(64, 7).Value = 190
$ws.Cells.Item(64, 8).Value = '• Lenovo 8GB DDR4 2400MHz ECC RDIMM Memory (4X70M09261) 
• Capacity: 8GB
• Memory Type: DDR4
• Weight: 0.02kg 
• Height: 3.98 mm 
• Width: 30.75 mm 
• Depth: 133.35 mm  '
$ws.Cells.Item(64, 8).WrapText = $true
$ws.Cells.Item(64, 9).Value = $brandVal
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Rows.Item(64).RowHeight = 116

# Row 65: 'Lenovo 16 GB DDR4 RAM Memory'
$ws.Cells.Item(65, 1).Value = 63
$ws.Cells.Item(65, 2).Value = $typeVal
$ws.Cells.Item(65, 3).Value = 'Lenovo 16 GB DDR4 RAM Memory'
$ws.Cells.Item(65, 4).Value = $noneVal
$ws.Cells.Item(65, 5).Value = $noneVal
$ws.Cells.Item(65, 6).Value = $noneVal
$ws.Cells.Item(65, 7).Value = 190
$ws.Cells.Item(65, 8).Value = '• Lenovo 16GB DDR4 2933MHz ECC SoDIMM Memory upgrades
   the standard memory capacity in select PCs, improving 
   performance and enhancing the PC''s ability to run more 
   programs simultaneously.
• 16GB DDR4 2933MHz ECC SoDIMM.
•  Compatible with select ThinkPad Mobile Workstations.
• Hardware Requirements: SoDIMM Slot.
• Weight: 15g '
$ws.Cells.Item(65, 8).WrapText = $true
$ws.Cells.Item(65, 9).Value = $brandVal
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Rows.Item(65).RowHeight = 174

# Update the saved selection/view to match the authored workbook
[void]$ws.Range("L1").Select()
